$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 481.87
$ws.Range("J17").Value = 481.87
$ws.Range("L17").Value = 1445.61
$ws.Range("N17").Value = -1781.61
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H64").Value = 2890.7273
$ws.Range("I64").Value = 2690
$ws.Range("J64").Value = 2935.3333
$ws.Range("K64").Value = 2690
$ws.Range("L64").Value = 2935.3333
$ws.Range("M64").Value = -2442
$ws.Range("N64").Value = -3431.3333
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H67").Value = 2890.7273
$ws.Range("I67").Value = 2690
$ws.Range("J67").Value = 2935.3333
$ws.Range("K67").Value = 2690
$ws.Range("L67").Value = 2935.3333
$ws.Range("M67").Value = -1832
$ws.Range("N67").Value = -4651.3333
$ws.Range("H69").Value = 4183
$ws.Range("I69").Value = 3515
$ws.Range("J69").Value = 4350
$ws.Range("K69").Value = 10545
$ws.Range("L69").Value = 13050
$ws.Range("M69").Value = -9671
$ws.Range("N69").Value = -14798
$ws.Range("H70").Value = 880.9524
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 875
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 2625
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -3165
$ws.Range("H72").Value = 4183
$ws.Range("I72").Value = 3515
$ws.Range("J72").Value = 4350
$ws.Range("K72").Value = 31635
$ws.Range("L72").Value = 39150
$ws.Range("M72").Value = -27267
$ws.Range("N72").Value = -47886
$ws.Range("H73").Value = 880.9524
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 875
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 2625
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -4497
$ws.Range("H137").Value = 924.56714
$ws.Range("I137").Value = 746.614
$ws.Range("J137").Value = 1938.9
$ws.Range("K137").Value = 2239.842
$ws.Range("L137").Value = 5816.700000000001
$ws.Range("M137").Value = 310.1579999999999
$ws.Range("N137").Value = -10916.7
$ws.Range("H138").Value = 18185738
$ws.Range("I138").Value = 31251248
$ws.Range("J138").Value = 7639.174
$ws.Range("K138").Value = 93753744
$ws.Range("L138").Value = 22917.522
$ws.Range("M138").Value = -93748604
$ws.Range("N138").Value = -33197.522
$ws.Range("H141").Value = 2886.4412
$ws.Range("I141").Value = 1549.8667
$ws.Range("J141").Value = 5501.478
$ws.Range("K141").Value = 4649.6001
$ws.Range("L141").Value = 16504.434
$ws.Range("M141").Value = 530.3999000000003
$ws.Range("N141").Value = -26864.434
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 861.55
$ws.Range("I74").Value = 773.5217
$ws.Range("J74").Value = 980.64703
$ws.Range("K74").Value = 773.5217
$ws.Range("L74").Value = 980.64703
$ws.Range("M74").Value = 100.4783
$ws.Range("N74").Value = -2728.64703
$ws.Range("H77").Value = 861.55
$ws.Range("I77").Value = 773.5217
$ws.Range("J77").Value = 980.64703
$ws.Range("K77").Value = 3867.6085
$ws.Range("L77").Value = 4903.23515
$ws.Range("M77").Value = 500.3914999999997
$ws.Range("N77").Value = -13639.23515
$ws.Range("H135").Value = 49425
$ws.Range("J135").Value = 49425
$ws.Range("L135").Value = 49425
$ws.Range("N135").Value = -59565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 947.1129
$ws.Range("I134").Value = 650.34
$ws.Range("K134").Value = 1951.02
$ws.Range("M134").Value = 583.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110.933334
$ws.Range("I7").Value = 110.2381
$ws.Range("J7").Value = 112.55556
$ws.Range("K7").Value = 110.2381
$ws.Range("L7").Value = 112.55556
$ws.Range("M7").Value = 2.761899999999997
$ws.Range("N7").Value = -338.55556
$ws.Range("H62").Value = 16718903
$ws.Range("I62").Value = 30396324
$ws.Range("K62").Value = 30396324
$ws.Range("M62").Value = -30395700
$ws.Range("H65").Value = 16718903
$ws.Range("I65").Value = 30396324
$ws.Range("K65").Value = 151981620
$ws.Range("M65").Value = -151978500
$ws.Range("H105").Value = 5435.125
$ws.Range("I105").Value = 8117.5
$ws.Range("J105").Value = 2752.75
$ws.Range("K105").Value = 8117.5
$ws.Range("L105").Value = 2752.75
$ws.Range("M105").Value = -6370.5
$ws.Range("N105").Value = -6246.75
$ws.Range("H134").Value = 861.26666
$ws.Range("I134").Value = 724.03125
$ws.Range("J134").Value = 1199.0769
$ws.Range("K134").Value = 2172.09375
$ws.Range("L134").Value = 3597.2307
$ws.Range("M134").Value = 362.90625
$ws.Range("N134").Value = -8667.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1085.5714
$ws.Range("I92").Value = 562.5
$ws.Range("J92").Value = 1294.8
$ws.Range("K92").Value = 1687.5
$ws.Range("L92").Value = 3884.4
$ws.Range("M92").Value = -439.5
$ws.Range("N92").Value = -6380.4
$ws.Range("H93").Value = 7556
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7556
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 22668
$ws.Range("N93").Value = -26412
$ws.Range("H94").Value = 6516.2
$ws.Range("J94").Value = 8020.25
$ws.Range("L94").Value = 24060.75
$ws.Range("N94").Value = -25412.75
$ws.Range("H95").Value = 7365.4
$ws.Range("J95").Value = 7365.4
$ws.Range("L95").Value = 22096.2
$ws.Range("N95").Value = -26214.2
$ws.Range("H96").Value = 6225
$ws.Range("J96").Value = 6225
$ws.Range("L96").Value = 18675
$ws.Range("N96").Value = -22793
$ws.Range("H97").Value = 458
$ws.Range("J97").Value = 454.5
$ws.Range("L97").Value = 1363.5
$ws.Range("N97").Value = -2355.5
$ws.Range("H98").Value = 324.33334
$ws.Range("I98").Value = 324.33334
$ws.Range("K98").Value = 973.0000200000001
$ws.Range("M98").Value = 524.9999799999999
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("H100").Value = 7666.6665
$ws.Range("J100").Value = 7666.6665
$ws.Range("L100").Value = 22999.9995
$ws.Range("N100").Value = -24621.9995
$ws.Range("H101").Value = 8530
$ws.Range("J101").Value = 8530
$ws.Range("L101").Value = 25590
$ws.Range("N101").Value = -30458
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868
$ws.Range("H103").Value = 2063.0417
$ws.Range("I103").Value = 499
$ws.Range("J103").Value = 2474.6316
$ws.Range("K103").Value = 1497
$ws.Range("L103").Value = 7423.8948
$ws.Range("M103").Value = -618
$ws.Range("N103").Value = -9181.8948
$ws.Range("H104").Value = 1089.1818
$ws.Range("I104").Value = 812
$ws.Range("J104").Value = 1193.125
$ws.Range("K104").Value = 2436
$ws.Range("L104").Value = 3579.375
$ws.Range("M104").Value = 185
$ws.Range("N104").Value = -8821.375
$ws.Range("H105").Value = 7822.6665
$ws.Range("J105").Value = 7822.6665
$ws.Range("L105").Value = 23467.9995
$ws.Range("N105").Value = -28709.9995
$ws.Range("H106").Value = 8900
$ws.Range("J106").Value = 8900
$ws.Range("L106").Value = 26700
$ws.Range("N106").Value = -28592
$ws.Range("H131").Value = 794.84906
$ws.Range("I131").Value = 333.86667
$ws.Range("J131").Value = 976.8158
$ws.Range("K131").Value = 1001.60001
$ws.Range("L131").Value = 2930.4474
$ws.Range("M131").Value = 4038.39999
$ws.Range("N131").Value = -13010.4474
$ws.Range("M93").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 626
$ws.Range("I107").Value = 516.75
$ws.Range("J107").Value = 844.5
$ws.Range("K107").Value = 516.75
$ws.Range("L107").Value = 844.5
$ws.Range("M107").Value = 1403.25
$ws.Range("N107").Value = -4684.5

Write-Host "Done applying Valefor_Profits updates"